$d = $word.ActiveDocument

function Get-ParaText($para) {
    return $para.Range.Text.TrimEnd([char]13, [char]10, [char]7)
}

# --- 1. Fill in the empty paragraph that immediately follows the
#        "Design Philosophy" heading. ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ((Get-ParaText $p) -eq "Design Philosophy") {
        $next = $d.Paragraphs.Item($i + 1)
        if ((Get-ParaText $next) -eq "") {
            $next.Range.Text = "The overall design philosophy consisted of using functions to make the code cleaner and more modular"
        }
        break
    }
}

# --- 2. Locate the two empty trailing paragraphs (right after the
#        "... documentation portion of the assignment." paragraph) and
#        turn them into a "References" heading + repository link. ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = Get-ParaText $p
    if ($t -like "*documentation portion of the assignment*") {
        $refHeading = $d.Paragraphs.Item($i + 1)
        $refLink = $d.Paragraphs.Item($i + 2)

        if ((Get-ParaText $refHeading) -eq "") {
            $refHeading.Range.Text = "References"
            $refHeading.Range.Font.Underline = 1
        }

        if ((Get-ParaText $refLink) -eq "") {
            $refLink.Range.Text = "Public Assignment Repository: https://github.com/Andy1184/COMP3100-Assignment-1"
        }
        break
    }
}
